$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.163221836090088
$ws.Range("B1").Value = 2.681084871292114
$ws.Range("C1").Value = 2.696366310119629
$ws.Range("D1").Value = 3.266540765762329
$ws.Range("E1").Value = 2.530049800872803
